# Updated cryptos list - apply price & volume changes from source diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.849.36"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "'1.735.65"
$ws.Range("E3").Value = "  -0.40%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'241.37"
$ws.Range("E5").Value = "  +4.70%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.5197"
$ws.Range("E7").Value = "  -0.87%  "

$ws.Range("D8").Value = "'0.2732"
$ws.Range("E8").Value = "  -0.71%  "

$ws.Range("D9").Value = "'0.06149"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").Value = "'1.740.47"
$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").Value = "'0.07165"
$ws.Range("E11").Value = "  +0.94%  "

$ws.Range("D12").Value = "'14.94"
$ws.Range("E12").Value = "  -1.75%  "

$ws.Range("D13").Value = "'0.6391"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").Value = "'4.602"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").Value = "'76.98"
$ws.Range("E15").Value = "  -0.37%  "

$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "'25.883.02"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").Value = "'11.72"
$ws.Range("E19").Value = "  +1.51%  "

$ws.Range("D20").Value = "'0.000006769"
$ws.Range("E20").Value = "  +1.27%  "

$ws.Range("D21").Value = "'1.962.77"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "'4.268"
$ws.Range("E22").Value = "  -0.56%  "

$ws.Range("D23").Value = "'8.597"
$ws.Range("E23").Value = "  -1.99%  "

$ws.Range("D24").Value = "'5.255"
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("D25").Value = "'137.52"
$ws.Range("E25").Value = "  -2.00%  "

$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("E28").Value = "  -2.00%  "

$ws.Range("D29").Value = "'104.89"
$ws.Range("E29").Value = "  +2.15%  "

$ws.Range("D30").Value = "'3.931"
$ws.Range("E30").Value = "  +5.26%  "

$ws.Range("D31").Value = "'0.08244"
$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("D32").Value = "'3.646"
$ws.Range("E32").Value = "  +2.19%  "

$ws.Range("D33").Value = "'0.04628"
$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("E34").Value = "  +1.15%  "

$ws.Range("D35").Value = "'0.9858"
$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("D36").Value = "'0.6176"
$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("D37").Value = "'2.687"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").Value = "'0.01594"
$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("E39").Value = "  -0.79%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").Value = "'100.17"
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").Value = "'0.3836"
$ws.Range("E42").Value = "  -1.06%  "

$ws.Range("D43").Value = "'0.7461"
$ws.Range("E43").Value = "  +1.55%  "

$ws.Range("D44").Value = "'4.992"
$ws.Range("E44").Value = "  -0.45%  "

$ws.Range("E45").Value = "  -0.48%  "

$ws.Range("D46").Value = "'6.228"
$ws.Range("E46").Value = "  -0.69%  "

$ws.Range("D47").Value = "'0.05234"
$ws.Range("E47").Value = "  -1.99%  "

$ws.Range("E48").Value = "  +2.05%  "

$ws.Range("D49").Value = "'30.51"
$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("D50").Value = "'7.562"
$ws.Range("E50").Value = "  -1.60%  "

$ws.Range("D51").Value = "'0.3398"
$ws.Range("E51").Value = "  -0.86%  "
